# 自动更新Excel文件 - 2025-12-20 23:12:10
# Decrement the "剩余" (remaining) value in column E by 1 for every data row,
# except row 36 (unchanged) and row 95 which is reset to a new cycle
# (remaining = 10, start date = 2025-12-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 99

for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 36) {
        continue
    }
    if ($r -eq 95) {
        $ws.Cells.Item($r, 5).Value2 = 10
        $ws.Cells.Item($r, 6).Value2 = 20251221
        continue
    }

    $cell = $ws.Cells.Item($r, 5)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current - 1
    }
}
